$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39 and 40 swap places (WhiteBITCoin <-> PolygonEcosystemToken)
# plus updated Price/Volume values for both.
$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.384'
$ws.Range("E39").Value = '  +11.68%  '

$ws.Range("B40").Value = 'WhiteBITCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '20.01'
$ws.Range("E40").Value = '  +1.97%  '

# --- Price / Volume(1h) updates for the remaining rows ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '76.413.39'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.033.09'
$ws.Range("E3").Value = '  +4.60%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '199.91'
$ws.Range("E5").Value = '  +1.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '628.51'
$ws.Range("E6").Value = '  +5.13%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.551'
$ws.Range("E8").Value = '  +0.61%  '
$ws.Range("E9").Value = '  +2.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.033.34'
$ws.Range("E10").Value = '  +4.70%  '
$ws.Range("E11").Value = '  +1.50%  '
$ws.Range("E12").Value = '  -0.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.08'
$ws.Range("E13").Value = '  +4.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.597.40'
$ws.Range("E14").Value = '  +4.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.19'
$ws.Range("E15").Value = '  +7.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '76.394.71'
$ws.Range("E16").Value = '  +1.00%  '
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.030.34'
$ws.Range("E18").Value = '  +4.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.42'
$ws.Range("E19").Value = '  +3.64%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.02'
$ws.Range("E20").Value = '  +3.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '372.95'
$ws.Range("E21").Value = '  +0.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.37'
$ws.Range("E22").Value = '  +2.07%  '
$ws.Range("E23").Value = '  -1.21%  '
$ws.Range("E24").Value = '  +4.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.20'
$ws.Range("E25").Value = '  +3.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.37'
$ws.Range("E27").Value = '  +4.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.89'
$ws.Range("E28").Value = '  +3.16%  '
$ws.Range("E29").Value = '  +0.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.30'
$ws.Range("E31").Value = '  +8.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.40'
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '510.87'
$ws.Range("E33").Value = '  +1.75%  '
$ws.Range("E34").Value = '  +7.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '20.74'
$ws.Range("E36").Value = '  +3.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '164.10'
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '194.22'
$ws.Range("E38").Value = '  +8.28%  '
$ws.Range("E41").Value = '  +0.94%  '
$ws.Range("E42").Value = '  -0.50%  '
$ws.Range("E43").Value = '  +0.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.06'
$ws.Range("E44").Value = '  +2.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.52'
$ws.Range("E45").Value = '  +6.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.25'
$ws.Range("E46").Value = '  +5.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.66'
$ws.Range("E47").Value = '  +0.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.716'
$ws.Range("E48").Value = '  +9.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.604'
$ws.Range("E49").Value = '  +6.25%  '
$ws.Range("E50").Value = '  +1.86%  '
$ws.Range("E51").Value = '  +4.60%  '
